# Auto-generated edit script: update profit/price figures across sheets
# per the authoritative diff (Pandaemonium_Profits.xlsx commit).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1700.8334
$ws.Range("I98").Value = 1700.8334
$ws.Range("K98").Value = 1700.8334
$ws.Range("M98").Value = -202.8334
$ws.Range("H101").Value = 2987.2104
$ws.Range("I101").Value = 1062.3334
$ws.Range("J101").Value = 3875.6155
$ws.Range("K101").Value = 3187.0002
$ws.Range("L101").Value = 11626.8465
$ws.Range("M101").Value = -1565.0002
$ws.Range("N101").Value = -14870.8465
$ws.Range("H122").Value = 1700.8334
$ws.Range("I122").Value = 1700.8334
$ws.Range("K122").Value = 5102.5002
$ws.Range("M122").Value = -2652.5002
$ws.Range("H129").Value = 1592.8572
$ws.Range("J129").Value = 2500
$ws.Range("L129").Value = 7500
$ws.Range("N129").Value = -17500
$ws.Range("H132").Value = 1406.1428
$ws.Range("I132").Value = 1232.2549
$ws.Range("J132").Value = 3179.8
$ws.Range("K132").Value = 3696.7647
$ws.Range("L132").Value = 9539.400000000001
$ws.Range("M132").Value = -1166.7647
$ws.Range("N132").Value = -14599.4
$ws.Range("H137").Value = 3112.9556
$ws.Range("I137").Value = 1555.9286
$ws.Range("J137").Value = 5677.4707
$ws.Range("K137").Value = 4667.7858
$ws.Range("L137").Value = 17032.4121
$ws.Range("M137").Value = -2117.7858
$ws.Range("N137").Value = -22132.4121
$ws.Range("H138").Value = 3842.8406
$ws.Range("I138").Value = 1670.3636
$ws.Range("J138").Value = 4254.8623
$ws.Range("K138").Value = 5011.0908
$ws.Range("L138").Value = 12764.5869
$ws.Range("M138").Value = 128.9092000000001
$ws.Range("N138").Value = -23044.5869

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7789.4375
$ws.Range("I32").Value = 5841.058
$ws.Range("J32").Value = 20011.092
$ws.Range("K32").Value = 5841.058
$ws.Range("L32").Value = 20011.092
$ws.Range("M32").Value = -5554.058
$ws.Range("N32").Value = -20585.092
$ws.Range("H45").Value = 2400
$ws.Range("I45").Value = 2000
$ws.Range("K45").Value = 2000
$ws.Range("M45").Value = -1623
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H122").Value = 15628708
$ws.Range("I122").Value = 3590.5
$ws.Range("K122").Value = 10771.5
$ws.Range("M122").Value = -8321.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H111").Value = 58966.332
$ws.Range("J111").Value = 58966.332
$ws.Range("L111").Value = 58966.332
$ws.Range("N111").Value = -67146.33199999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2479.66
$ws.Range("I31").Value = 1769.9756
$ws.Range("J31").Value = 5712.6665
$ws.Range("K31").Value = 1769.9756
$ws.Range("L31").Value = 5712.6665
$ws.Range("M31").Value = -1474.9756
$ws.Range("N31").Value = -6302.6665
$ws.Range("H34").Value = 2479.66
$ws.Range("I34").Value = 1769.9756
$ws.Range("J34").Value = 5712.6665
$ws.Range("K34").Value = 1769.9756
$ws.Range("L34").Value = 5712.6665
$ws.Range("M34").Value = -1567.9756
$ws.Range("N34").Value = -6116.6665
$ws.Range("H58").Value = 2069311
$ws.Range("I58").Value = 3137203
$ws.Range("J58").Value = 4719.933
$ws.Range("K58").Value = 3137203
$ws.Range("L58").Value = 4719.933
$ws.Range("M58").Value = -3137000
$ws.Range("N58").Value = -5125.933
$ws.Range("H74").Value = 36314
$ws.Range("J74").Value = 36314
$ws.Range("L74").Value = 36314
$ws.Range("N74").Value = -38062
$ws.Range("H77").Value = 36314
$ws.Range("J77").Value = 36314
$ws.Range("L77").Value = 108942
$ws.Range("N77").Value = -117678
$ws.Range("H99").Value = 2122.4285
$ws.Range("I99").Value = 2860
$ws.Range("J99").Value = 1712.6666
$ws.Range("K99").Value = 2860
$ws.Range("L99").Value = 1712.6666
$ws.Range("M99").Value = -1362
$ws.Range("N99").Value = -4708.6666
$ws.Range("H107").Value = 632.7
$ws.Range("I107").Value = 543.4286
$ws.Range("J107").Value = 841
$ws.Range("K107").Value = 543.4286
$ws.Range("L107").Value = 841
$ws.Range("M107").Value = 1376.5714
$ws.Range("N107").Value = -4681
$ws.Range("H126").Value = 2122.4285
$ws.Range("I126").Value = 2860
$ws.Range("J126").Value = 1712.6666
$ws.Range("K126").Value = 8580
$ws.Range("L126").Value = 5137.9998
$ws.Range("M126").Value = -6110
$ws.Range("N126").Value = -10077.9998
$ws.Range("H136").Value = 2069311
$ws.Range("I136").Value = 3137203
$ws.Range("J136").Value = 4719.933
$ws.Range("K136").Value = 9411609
$ws.Range("L136").Value = 14159.799
$ws.Range("M136").Value = -9409059
$ws.Range("N136").Value = -19259.799

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 95
$ws.Range("I8").Value = 95
$ws.Range("K8").Value = 285
$ws.Range("M8").Value = -146

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 4500
$ws.Range("I22").Value = 4500
$ws.Range("K22").Value = 4500
$ws.Range("M22").Value = -3971
$ws.Range("H47").Value = 17160.572
$ws.Range("J47").Value = 17160.572
$ws.Range("L47").Value = 17160.572
$ws.Range("N47").Value = -18296.572
$ws.Range("H55").Value = 14666.667
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 14666.667
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 14666.667
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -15320.667
$ws.Range("H102").Value = 3254.775
$ws.Range("I102").Value = 3322.44
$ws.Range("J102").Value = 3142
$ws.Range("K102").Value = 3322.44
$ws.Range("L102").Value = 3142
$ws.Range("M102").Value = -1700.44
$ws.Range("N102").Value = -6386
$ws.Range("H112").Value = 79800
$ws.Range("J112").Value = 79800
$ws.Range("L112").Value = 79800
$ws.Range("N112").Value = -82016
$ws.Range("H122").Value = 15000
$ws.Range("I122").Value = 50000
$ws.Range("J122").Value = 6250
$ws.Range("K122").Value = 150000
$ws.Range("L122").Value = 18750
$ws.Range("M122").Value = -147550
$ws.Range("N122").Value = -23650
$ws.Range("H126").Value = 2980.4092
$ws.Range("I126").Value = 1875
$ws.Range("J126").Value = 4577.1113
$ws.Range("K126").Value = 5625
$ws.Range("L126").Value = 13731.3339
$ws.Range("M126").Value = -3155
$ws.Range("N126").Value = -18671.3339

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 37048.6
$ws.Range("J62").Value = 37048.6
$ws.Range("L62").Value = 37048.6
$ws.Range("N62").Value = -38296.6
$ws.Range("H64").Value = 34074.75
$ws.Range("J64").Value = 34074.75
$ws.Range("L64").Value = 34074.75
$ws.Range("N64").Value = -34524.75
$ws.Range("H65").Value = 37048.6
$ws.Range("J65").Value = 37048.6
$ws.Range("L65").Value = 111145.8
$ws.Range("N65").Value = -117385.8
$ws.Range("H67").Value = 34074.75
$ws.Range("J67").Value = 34074.75
$ws.Range("L67").Value = 34074.75
$ws.Range("N67").Value = -35634.75
$ws.Range("H68").Value = 2164.2856
$ws.Range("I68").Value = 1830
$ws.Range("K68").Value = 1830
$ws.Range("M68").Value = -1081
$ws.Range("H71").Value = 2164.2856
$ws.Range("I71").Value = 1830
$ws.Range("K71").Value = 9150
$ws.Range("M71").Value = -5406
$ws.Range("H76").Value = 27332.111
$ws.Range("I76").Value = 10261
$ws.Range("J76").Value = 29466
$ws.Range("K76").Value = 10261
$ws.Range("L76").Value = 29466
$ws.Range("M76").Value = -9923
$ws.Range("N76").Value = -30142
$ws.Range("H79").Value = 27332.111
$ws.Range("I79").Value = 10261
$ws.Range("J79").Value = 29466
$ws.Range("K79").Value = 10261
$ws.Range("L79").Value = 29466
$ws.Range("M79").Value = -9091
$ws.Range("N79").Value = -31806
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H122").Value = 5876.976
$ws.Range("I122").Value = 4609
$ws.Range("J122").Value = 7937.4375
$ws.Range("K122").Value = 13827
$ws.Range("L122").Value = 23812.3125
$ws.Range("M122").Value = -11377
$ws.Range("N122").Value = -28712.3125

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 11250
$ws.Range("J41").Value = 11250
$ws.Range("L41").Value = 11250
$ws.Range("N41").Value = -12030
